$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 216, pushing existing rows 216-230 down to 218-232.
$ws.Rows.Item(216).Insert()
$ws.Rows.Item(216).Insert()

# Populate the new row 216 (Brócoli, Primera, new date 44516)
$ws.Cells.Item(216, 1).Value = 4
$ws.Cells.Item(216, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216, 3).Value = "Los Lagos"
$ws.Cells.Item(216, 4).Value = 44516
$ws.Cells.Item(216, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 5).Value = 10
$ws.Cells.Item(216, 6).Value = 100112023
$ws.Cells.Item(216, 7).Value = "Brócoli"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 700
$ws.Cells.Item(216, 11).Value = 1200
$ws.Cells.Item(216, 12).Value = 1200
$ws.Cells.Item(216, 13).Value = 1200
$ws.Cells.Item(216, 14).Value = "$/unidad"
$ws.Cells.Item(216, 15).Value = "Región Metropolitana"
$ws.Cells.Item(216, 16).Value = 1200
$ws.Cells.Item(216, 17).Value = 1
$ws.Cells.Item(216, 18).Value = "Hortaliza"

# Populate the new row 217 (Brócoli, Segunda, new date 44516)
$ws.Cells.Item(217, 1).Value = 4
$ws.Cells.Item(217, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(217, 3).Value = "Los Lagos"
$ws.Cells.Item(217, 4).Value = 44516
$ws.Cells.Item(217, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 5).Value = 10
$ws.Cells.Item(217, 6).Value = 100112023
$ws.Cells.Item(217, 7).Value = "Brócoli"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Segunda"
$ws.Cells.Item(217, 10).Value = 700
$ws.Cells.Item(217, 11).Value = 1000
$ws.Cells.Item(217, 12).Value = 1000
$ws.Cells.Item(217, 13).Value = 1000
$ws.Cells.Item(217, 14).Value = "$/unidad"
$ws.Cells.Item(217, 15).Value = "Región Metropolitana"
$ws.Cells.Item(217, 16).Value = 1000
$ws.Cells.Item(217, 17).Value = 1
$ws.Cells.Item(217, 18).Value = "Hortaliza"
